$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current totals row (row 15), shifting the
# totals row down to row 16, and copying the previous data row's
# formatting into the freshly inserted row.
$ws.Range("A15:F15").Insert()
$ws.Range("A14:F14").Copy()
$ws.Range("A15:F15").PasteSpecial(-4122) # xlPasteFormats

# Populate the new data row (row 15) with the new time entry.
$ws.Cells.Item(15, 1).Value = 45280
$ws.Cells.Item(15, 2).Value = 0.583333333333333
$ws.Cells.Item(15, 3).Value = 0.916666666666667
$ws.Cells.Item(15, 4).Formula = "=(C15<B15)+C15-B15"
$ws.Cells.Item(15, 5).Value = 10
$ws.Cells.Item(15, 6).Formula = "=(D15*24)*E15"

# PasteSpecial(xlPasteFormats) mirrors column D's style onto F; restore F15
# to the same style used by the other data rows' Bill column (style id 9,
# same numFmt/border as E15, which already pasted correctly).
$ws.Range("E15").Copy()
$ws.Range("F15").PasteSpecial(-4122) # xlPasteFormats

# Fix up the totals row (now row 16) so its SUM ranges include the new row.
$ws.Cells.Item(16, 4).Formula = "=SUM(D2:D15)"
$ws.Cells.Item(16, 6).Formula = "=SUM(F2:F15)"

$wb.Application.Calculate()

$ws.Range("F17").Select()
